$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.24418859991079245
$ws.Range("A2").Value = -0.0059999999812774263
$ws.Range("A3").Value = -0.0039999999815574228
$ws.Range("A4").Value = -0.0079999999680016032
$ws.Range("A5").Value = -0.0029999999809149358
$ws.Range("A6").Value = 0.014822428103748564
$ws.Range("A7").Value = -0.0099999999560225206
$ws.Range("A8").Value = -0.0099999999562463415
$ws.Range("A9").Value = -0.0019999999806437074
$ws.Range("A10").Value = -0.029306705808506806
$ws.Range("A11").Value = -0.0029999999791101573
$ws.Range("A12").Value = 0.059237003992740789
$ws.Range("A13").Value = -0.0034999999779463664
$ws.Range("A14").Value = -0.0079999999645803399
$ws.Range("A15").Value = 0.011674708627449348
$ws.Range("A16").Value = -0.0019999999819781955
$ws.Range("A17").Value = -0.0019999999815381031
$ws.Range("A18").Value = -0.0039999999755169213
$ws.Range("A19").Value = -0.0039999999870863334
$ws.Range("A20").Value = -0.0039999999861457525
$ws.Range("A21").Value = -0.003999999985999203
$ws.Range("A22").Value = -0.0039999999858562063
$ws.Range("A23").Value = -0.077106449201973071
$ws.Range("A24").Value = -0.019999999929861012
$ws.Range("A25").Value = -0.019999999928890233
$ws.Range("A26").Value = -0.00249999997715733
$ws.Range("A27").Value = -0.0024999999758756886
$ws.Range("A28").Value = -0.0019999999717175143
$ws.Range("A29").Value = -0.0069999999526570988
$ws.Range("A30").Value = -0.059999999794574599
$ws.Range("A31").Value = -0.0069999999490804043
$ws.Range("A32").Value = -0.0099999999399909001
$ws.Range("A33").Value = -0.0039999999575659473
